$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.7943689804753
$ws.Range("D2").Value = 5.776676400080993
$ws.Range("E2").Value = 15.90878481530115
$ws.Range("F2").Value = 29.69623103066452
$ws.Range("G2").Value = 38.31266637593834
$ws.Range("H2").Value = 16.50102973411032
$ws.Range("I2").Value = 29.1126961738474
$ws.Range("K2").Value = 9.625368628760203
$ws.Range("L2").Value = 9.445300538166702
$ws.Range("M2").Value = 14.24889270605433

$ws.Range("B3").Value = 12.70648831570122
$ws.Range("D3").Value = 5.767609287115135
$ws.Range("E3").Value = 15.915714239621
$ws.Range("F3").Value = 29.64031569518267
$ws.Range("G3").Value = 38.18820906881495
$ws.Range("H3").Value = 16.53087260967802
$ws.Range("I3").Value = 29.20698031104261
$ws.Range("K3").Value = 9.221443738200657
$ws.Range("L3").Value = 9.433834817097527
$ws.Range("M3").Value = 14.23412268442984

$ws.Range("B4").Value = 12.65461005694017
$ws.Range("D4").Value = 5.761935856030408
$ws.Range("E4").Value = 15.92042833643783
$ws.Range("F4").Value = 29.61401989496536
$ws.Range("G4").Value = 38.12359773278132
$ws.Range("H4").Value = 16.55228940541043
$ws.Range("I4").Value = 29.26940430333186
$ws.Range("K4").Value = 8.962009493156513
$ws.Range("L4").Value = 9.428403566259053
$ws.Range("M4").Value = 14.22718692161551

$ws.Range("B5").Value = 12.63401162208268
$ws.Range("D5").Value = 5.759597547638085
$ws.Range("E5").Value = 15.92246507631469
$ws.Range("F5").Value = 29.60533096756712
$ws.Range("G5").Value = 38.10025298203777
$ws.Range("H5").Value = 16.56179347173277
$ws.Range("I5").Value = 29.29598183880119
$ws.Range("K5").Value = 8.853495034143169
$ws.Range("L5").Value = 9.426596786864835
$ws.Range("M5").Value = 14.22489988826488

$ws.Range("B6").Value = 12.63062456839261
$ws.Range("D6").Value = 5.759207693148831
$ws.Range("E6").Value = 15.92281026940321
$ws.Range("F6").Value = 29.6040107276624
$ws.Range("G6").Value = 38.09655727626648
$ws.Range("H6").Value = 16.56341847702274
$ws.Range("I6").Value = 29.30046380759384
$ws.Range("K6").Value = 8.835310163010602
$ws.Range("L6").Value = 9.42632137670293
$ws.Range("M6").Value = 14.22455278447781

$ws.Range("B7").Value = 12.65433003852984
$ws.Range("D7").Value = 5.76190442696256
$ws.Range("E7").Value = 15.92045533585719
$ws.Range("F7").Value = 29.61389450018775
$ws.Range("G7").Value = 38.12327079334669
$ws.Range("H7").Value = 16.55241443814464
$ws.Range("I7").Value = 29.26975812535949
$ws.Range("K7").Value = 8.960557223852739
$ws.Range("L7").Value = 9.428377551088856
$ws.Range("M7").Value = 14.22715389038403

$ws.Range("B8").Value = 12.76364758418016
$ws.Range("D8").Value = 5.773572313357357
$ws.Range("E8").Value = 15.91107884259672
$ws.Range("F8").Value = 29.67528767152083
$ws.Range("G8").Value = 38.26731412627783
$ws.Range("H8").Value = 16.51067686255677
$ws.Range("I8").Value = 29.14426410304652
$ws.Range("K8").Value = 9.48851017916517
$ws.Range("L8").Value = 9.441014512641024
$ws.Range("M8").Value = 14.24335913602519

$ws.Range("B9").Value = 12.99365126737829
$ws.Range("D9").Value = 5.795595375493309
$ws.Range("E9").Value = 15.89632889110162
$ws.Range("F9").Value = 29.85909002134677
$ws.Range("G9").Value = 38.64255552727487
$ws.Range("H9").Value = 16.45342367521233
$ws.Range("I9").Value = 28.93416778940192
$ws.Range("K9").Value = 10.43026966349922
$ws.Range("L9").Value = 9.478464277505616
$ws.Range("M9").Value = 14.29192235362139

$ws.Range("B10").Value = 13.17094192373882
$ws.Range("D10").Value = 5.811233869591749
$ws.Range("E10").Value = 15.88769954985891
$ws.Range("F10").Value = 30.03217848875862
$ws.Range("G10").Value = 38.97332409559384
$ws.Range("H10").Value = 16.42641628506745
$ws.Range("I10").Value = 28.80179346976756
$ws.Range("K10").Value = 11.06193582739408
$ws.Range("L10").Value = 9.513559252122281
$ws.Range("M10").Value = 14.33763128015505

$ws.Range("B11").Value = 13.25312824130732
$ws.Range("D11").Value = 5.818226636444771
$ws.Range("E11").Value = 15.88425115537745
$ws.Range("F11").Value = 30.11901052980644
$ws.Range("G11").Value = 39.13536274916107
$ws.Range("H11").Value = 16.41740972620707
$ws.Range("I11").Value = 28.74635499551347
$ws.Range("K11").Value = 11.33570756252817
$ws.Range("L11").Value = 9.531135489780269
$ws.Range("M11").Value = 14.36055439607462

$ws.Range("B12").Value = 13.28444751987
$ws.Range("D12").Value = 5.820856826113037
$ws.Range("E12").Value = 15.88301378786595
$ws.Range("F12").Value = 30.15303759284533
$ws.Range("G12").Value = 39.19834732503866
$ws.Range("H12").Value = 16.41447126517808
$ws.Range("I12").Value = 28.7260499521951
$ws.Range("K12").Value = 11.43739215241315
$ws.Range("L12").Value = 9.538019444109823
$ws.Range("M12").Value = 14.3695364377245

$ws.Range("B13").Value = 13.27769396713852
$ws.Range("D13").Value = 5.820291166273475
$ws.Range("E13").Value = 15.88327723442287
$ws.Range("F13").Value = 30.14565861207324
$ws.Range("G13").Value = 39.18471090825796
$ws.Range("H13").Value = 16.41508311008636
$ws.Range("I13").Value = 28.73039237730946
$ws.Range("K13").Value = 11.41558140805554
$ws.Range("L13").Value = 9.536526769682197
$ws.Range("M13").Value = 14.36758866063993

$ws.Range("B14").Value = 13.25570108379398
$ws.Range("D14").Value = 5.81844338303597
$ws.Range("E14").Value = 15.88414798510079
$ws.Range("F14").Value = 30.12178710692485
$ws.Range("G14").Value = 39.1405122193765
$ws.Range("H14").Value = 16.41715851245682
$ws.Range("I14").Value = 28.74467068671715
$ws.Range("K14").Value = 11.34411321053124
$ws.Range("L14").Value = 9.531697283962414
$ws.Range("M14").Value = 14.36128733843328

$ws.Range("B15").Value = 13.24225476759981
$ws.Range("D15").Value = 5.817309230557021
$ws.Range("E15").Value = 15.8846902568562
$ws.Range("F15").Value = 30.10731374726536
$ws.Range("G15").Value = 39.11364950833958
$ws.Range("H15").Value = 16.41849125527375
$ws.Range("I15").Value = 28.75350623221345
$ws.Range("K15").Value = 11.3000772097824
$ws.Range("L15").Value = 9.528768696930639
$ws.Range("M15").Value = 14.35746672328895

$ws.Range("B16").Value = 13.16559973549338
$ws.Range("D16").Value = 5.810774402994011
$ws.Range("E16").Value = 15.88793449910201
$ws.Range("F16").Value = 30.02666516056385
$ws.Range("G16").Value = 38.96296399430147
$ws.Range("H16").Value = 16.42707092496588
$ws.Range("I16").Value = 28.80551276526359
$ws.Range("K16").Value = 11.04376791549948
$ws.Range("L16").Value = 9.512442734869957
$ws.Range("M16").Value = 14.33617566217627

$ws.Range("B17").Value = 13.11895091458334
$ws.Range("D17").Value = 5.806734166928512
$ws.Range("E17").Value = 15.89004683688312
$ws.Range("F17").Value = 29.97925044599956
$ws.Range("G17").Value = 38.87345957791226
$ws.Range("H17").Value = 16.4331746337274
$ws.Range("I17").Value = 28.83864192038634
$ws.Range("K17").Value = 10.88302667310848
$ws.Range("L17").Value = 9.502837581966295
$ws.Range("M17").Value = 14.32365653049784

$ws.Range("B18").Value = 13.09226575410615
$ws.Range("D18").Value = 5.804398938801604
$ws.Range("E18").Value = 15.8913067166475
$ws.Range("F18").Value = 29.95274174094513
$ws.Range("G18").Value = 38.82307066826498
$ws.Range("H18").Value = 16.43699395972916
$ws.Range("I18").Value = 28.8581467182575
$ws.Range("K18").Value = 10.78929633655241
$ws.Range("L18").Value = 9.49746484537561
$ws.Range("M18").Value = 14.316656632906

$ws.Range("B19").Value = 13.08325641412861
$ws.Range("D19").Value = 5.803606325935101
$ws.Range("E19").Value = 15.89174100967517
$ws.Range("F19").Value = 29.94389791321642
$ws.Range("G19").Value = 38.80619849829865
$ws.Range("H19").Value = 16.43834010560727
$ws.Range("I19").Value = 28.86482792661122
$ws.Range("K19").Value = 10.75734281604982
$ws.Range("L19").Value = 9.495671919509519
$ws.Range("M19").Value = 14.31432121070635

$ws.Range("B20").Value = 13.12390182134091
$ws.Range("D20").Value = 5.807165438668656
$ws.Range("E20").Value = 15.88981732707622
$ws.Range("F20").Value = 29.98421898247389
$ws.Range("G20").Value = 38.88287476979021
$ws.Range("H20").Value = 16.43249293625071
$ws.Range("I20").Value = 28.8350687066733
$ws.Range("K20").Value = 10.9002702068923
$ws.Range("L20").Value = 9.503844371067387
$ws.Range("M20").Value = 14.32496846723012

$ws.Range("B21").Value = 13.26215576780812
$ws.Range("D21").Value = 5.818986608403884
$ws.Range("E21").Value = 15.88389036756329
$ws.Range("F21").Value = 30.12876780806568
$ws.Range("G21").Value = 39.15345070289127
$ws.Range("H21").Value = 16.41653609937273
$ws.Range("I21").Value = 28.74045811237447
$ws.Range("K21").Value = 11.36515932014055
$ws.Range("L21").Value = 9.533109656358
$ws.Range("M21").Value = 14.36313004435406

$ws.Range("B22").Value = 13.3536504382712
$ws.Range("D22").Value = 5.826608386306305
$ws.Range("E22").Value = 15.88041575565333
$ws.Range("F22").Value = 30.22990681094536
$ws.Range("G22").Value = 39.33973317548246
$ws.Range("H22").Value = 16.40885944214501
$ws.Range("I22").Value = 28.6826370704105
$ws.Range("K22").Value = 11.65739853224994
$ws.Range("L22").Value = 9.553564632609014
$ws.Range("M22").Value = 14.38982615673963

$ws.Range("B23").Value = 13.30472180762096
$ws.Range("D23").Value = 5.822550145027508
$ws.Range("E23").Value = 15.88223376056281
$ws.Range("F23").Value = 30.17532338747879
$ws.Range("G23").Value = 39.23946044863685
$ws.Range("H23").Value = 16.412704661636
$ws.Range("I23").Value = 28.71312976117867
$ws.Range("K23").Value = 11.50249560160272
$ws.Range("L23").Value = 9.542527092693398
$ws.Range("M23").Value = 14.37541895660892

$ws.Range("B24").Value = 13.12166309760848
$ws.Range("D24").Value = 5.806970499423689
$ws.Range("E24").Value = 15.88992094685896
$ws.Range("F24").Value = 29.98197036750224
$ws.Range("G24").Value = 38.87861483057622
$ws.Range("H24").Value = 16.43280016535483
$ws.Range("I24").Value = 28.83668272872097
$ws.Range("K24").Value = 10.89247850812318
$ws.Range("L24").Value = 9.503388736101112
$ws.Range("M24").Value = 14.32437472502653

$ws.Range("B25").Value = 12.92987951951035
$ws.Range("D25").Value = 5.789731343306193
$ws.Range("E25").Value = 15.89993078633872
$ws.Range("F25").Value = 29.80263259542574
$ws.Range("G25").Value = 38.53124770101194
$ws.Range("H25").Value = 16.46627204195483
$ws.Range("I25").Value = 28.98714671129626
$ws.Range("K25").Value = 10.18588059988359
$ws.Range("L25").Value = 9.466990566841556
$ws.Range("M25").Value = 14.2770076124564
